$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# StartDate (C) / EndDate (D) change for every data row (2-9)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = "10/02/2019"
    $ws.Cells.Item($r, 4).Value = "12/02/2019"
}

# Adults (E) becomes a text value instead of a number:
#   rows 2-5 -> "2"   rows 6-9 -> "3"
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 5).Value = "2"
}
for ($r = 6; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = "3"
}

# Kids (F) / AgeKid1 (G) drop back to 0 for rows 6-9.
# These columns are formatted as Text ("@"), so a plain .Value assignment
# of a number would be stored as a literal string. Flip the number format
# to General for the write, then restore it to Text to keep the original
# cell style/format intact (matches the source file's odd-but-intentional
# "numbers stored in text-formatted cells" layout).
for ($r = 6; $r -le 9; $r++) {
    $f = $ws.Cells.Item($r, 6)
    $f.NumberFormat = "General"
    $f.Value = 0
    $f.NumberFormat = "@"

    $g = $ws.Cells.Item($r, 7)
    $g.NumberFormat = "General"
    $g.Value = 0
    $g.NumberFormat = "@"
}

# Move the view/selection the way it ended up in the saved workbook
$null = $ws.Range("A1").Select()
$null = $ws.Range("E10").Select()
